$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 3251687.2
$ws.Range("J17").Value = 3313036.8
$ws.Range("L17").Value = 9939110.399999999
$ws.Range("N17").Value = -9939446.399999999
$ws.Range("H41").Value = 1067.3334
$ws.Range("J41").Value = 1301
$ws.Range("L41").Value = 1301
$ws.Range("N41").Value = -2181
$ws.Range("H106").Value = 70177970
$ws.Range("I106").Value = 47622264
$ws.Range("J106").Value = 83335460
$ws.Range("K106").Value = 47622264
$ws.Range("L106").Value = 83335460
$ws.Range("M106").Value = -47621633
$ws.Range("N106").Value = -83336722
$ws.Range("H112").Value = 1957.375
$ws.Range("J112").Value = 2426.5
$ws.Range("L112").Value = 7279.5
$ws.Range("N112").Value = -9495.5
$ws.Range("H129").Value = 875.9
$ws.Range("J129").Value = 968.17285
$ws.Range("L129").Value = 2904.51855
$ws.Range("N129").Value = -12904.51855
$ws.Range("H138").Value = 2062.723
$ws.Range("I138").Value = 709.8913
$ws.Range("J138").Value = 3744.6216
$ws.Range("K138").Value = 2129.6739
$ws.Range("L138").Value = 11233.8648
$ws.Range("M138").Value = 3010.3261
$ws.Range("N138").Value = -21513.8648

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 2431
$ws.Range("I2").Value = 1004.25
$ws.Range("J2").Value = 4333.3335
$ws.Range("K2").Value = 1004.25
$ws.Range("L2").Value = 4333.3335
$ws.Range("M2").Value = -891.25
$ws.Range("N2").Value = -4559.3335
$ws.Range("H63").Value = 62502300
$ws.Range("H66").Value = 62502300
$ws.Range("H76").Value = 29975
$ws.Range("J76").Value = 29975
$ws.Range("L76").Value = 29975
$ws.Range("N76").Value = -30651
$ws.Range("H79").Value = 29975
$ws.Range("J79").Value = 29975
$ws.Range("L79").Value = 29975
$ws.Range("N79").Value = -32315
$ws.Range("H116").Value = 2431
$ws.Range("I116").Value = 1004.25
$ws.Range("J116").Value = 4333.3335
$ws.Range("K116").Value = 1004.25
$ws.Range("L116").Value = 4333.3335
$ws.Range("M116").Value = 1289.75
$ws.Range("N116").Value = -8921.333500000001

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 2431
$ws.Range("I3").Value = 1004.25
$ws.Range("J3").Value = 4333.3335
$ws.Range("K3").Value = 1004.25
$ws.Range("L3").Value = 4333.3335
$ws.Range("M3").Value = -890.25
$ws.Range("N3").Value = -4561.3335
$ws.Range("H105").Value = 3477.6667
$ws.Range("I105").Value = 2500
$ws.Range("J105").Value = 4176
$ws.Range("K105").Value = 2500
$ws.Range("L105").Value = 4176
$ws.Range("M105").Value = -753
$ws.Range("N105").Value = -7670
$ws.Range("H134").Value = 2688.3103
$ws.Range("I134").Value = 2781.923
$ws.Range("K134").Value = 8345.769
$ws.Range("M134").Value = -5810.769

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 2748462
$ws.Range("I16").Value = 5129153
$ws.Range("J16").Value = 1511
$ws.Range("K16").Value = 5129153
$ws.Range("L16").Value = 1511
$ws.Range("M16").Value = -5128866
$ws.Range("N16").Value = -2085
$ws.Range("H31").Value = 3450.383
$ws.Range("I31").Value = 1748.2222
$ws.Range("J31").Value = 5748.3
$ws.Range("K31").Value = 1748.2222
$ws.Range("L31").Value = 5748.3
$ws.Range("M31").Value = -1453.2222
$ws.Range("N31").Value = -6338.3
$ws.Range("H34").Value = 3450.383
$ws.Range("I34").Value = 1748.2222
$ws.Range("J34").Value = 5748.3
$ws.Range("K34").Value = 1748.2222
$ws.Range("L34").Value = 5748.3
$ws.Range("M34").Value = -1546.2222
$ws.Range("N34").Value = -6152.3
$ws.Range("H113").Value = 2748462
$ws.Range("I113").Value = 5129153
$ws.Range("J113").Value = 1511
$ws.Range("K113").Value = 5129153
$ws.Range("L113").Value = 1511
$ws.Range("M113").Value = -5126983
$ws.Range("N113").Value = -5851
$ws.Range("H132").Value = 2859.1155
$ws.Range("I132").Value = 2632.4211
$ws.Range("K132").Value = 7897.263300000001
$ws.Range("M132").Value = -5367.263300000001
$ws.Range("H134").Value = 2939.25
$ws.Range("I134").Value = 2859.1428
$ws.Range("J134").Value = 3500
$ws.Range("K134").Value = 8577.428400000001
$ws.Range("L134").Value = 10500
$ws.Range("M134").Value = -6042.428400000001
$ws.Range("N134").Value = -15570

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 273408.12
$ws.Range("I5").Value = 756.1667
$ws.Range("J5").Value = 600590.5
$ws.Range("K5").Value = 2268.5001
$ws.Range("L5").Value = 1801771.5
$ws.Range("M5").Value = -2156.5001
$ws.Range("N5").Value = -1801995.5
$ws.Range("H12").Value = 3846276.8
$ws.Range("J12").Value = 132.35294
$ws.Range("L12").Value = 397.05882
$ws.Range("N12").Value = -743.05882
$ws.Range("H17").Value = 798.4
$ws.Range("J17").Value = 798.4
$ws.Range("L17").Value = 2395.2
$ws.Range("N17").Value = -2733.2
$ws.Range("H31").Value = 363.33334
$ws.Range("I31").Value = 363.33334
$ws.Range("K31").Value = 1090.00002
$ws.Range("M31").Value = -802.0000199999999
$ws.Range("H44").Value = 747.65216
$ws.Range("I44").Value = 461.3846
$ws.Range("J44").Value = 1119.8
$ws.Range("K44").Value = 1384.1538
$ws.Range("L44").Value = 3359.4
$ws.Range("M44").Value = -986.1538
$ws.Range("N44").Value = -4155.4
$ws.Range("H55").Value = 2343.6875
$ws.Range("J55").Value = 2343.6875
$ws.Range("L55").Value = 7031.0625
$ws.Range("N55").Value = -7385.0625
$ws.Range("H116").Value = 2633
$ws.Range("I116").Value = 0
$ws.Range("J116").Value = 2633
$ws.Range("K116").Value = 0
$ws.Range("L116").Value = 7899
$ws.Range("N116").Value = -14783
$ws.Range("M116").ClearContents()
$ws.Range("H132").Value = 1141.0769
$ws.Range("I132").Value = 723.4
$ws.Range("J132").Value = 2533.3333
$ws.Range("K132").Value = 6510.599999999999
$ws.Range("L132").Value = 22799.9997
$ws.Range("M132").Value = -3980.599999999999
$ws.Range("N132").Value = -27859.9997
$ws.Range("H135").Value = 273408.12
$ws.Range("I135").Value = 756.1667
$ws.Range("J135").Value = 600590.5
$ws.Range("K135").Value = 6805.5003
$ws.Range("L135").Value = 5405314.5
$ws.Range("M135").Value = -4270.5003
$ws.Range("N135").Value = -5410384.5
$ws.Range("H140").Value = 2823.4167
$ws.Range("I140").Value = 2171
$ws.Range("J140").Value = 10000
$ws.Range("K140").Value = 6513
$ws.Range("L140").Value = 30000
$ws.Range("M140").Value = -1333
$ws.Range("N140").Value = -40360

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("I74").Value = 40000
$ws.Range("K74").Value = 40000
$ws.Range("M74").Value = -39064
$ws.Range("I77").Value = 40000
$ws.Range("K77").Value = 120000
$ws.Range("M77").Value = -115320
$ws.Range("H80").Value = 2593.0356
$ws.Range("I80").Value = 2460.5
$ws.Range("K80").Value = 2460.5
$ws.Range("M80").Value = -1462.5
$ws.Range("H83").Value = 2593.0356
$ws.Range("I83").Value = 2460.5
$ws.Range("K83").Value = 12302.5
$ws.Range("M83").Value = -7310.5
$ws.Range("H132").Value = 3545.9167
$ws.Range("I132").Value = 3611.625
$ws.Range("J132").Value = 3414.5
$ws.Range("K132").Value = 10834.875
$ws.Range("L132").Value = 10243.5
$ws.Range("M132").Value = -8304.875
$ws.Range("N132").Value = -15303.5

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H14").Value = 26526.125
$ws.Range("I14").Value = 17700
$ws.Range("J14").Value = 35352.25
$ws.Range("K14").Value = 17700
$ws.Range("L14").Value = 35352.25
$ws.Range("M14").Value = -17532
$ws.Range("N14").Value = -35688.25
$ws.Range("H32").Value = 0
$ws.Range("J32").Value = 0
$ws.Range("L32").Value = 0
$ws.Range("N32").ClearContents()
$ws.Range("H86").Value = 17630
$ws.Range("J86").Value = 17630
$ws.Range("L86").Value = 17630
$ws.Range("N86").Value = -19876
$ws.Range("H89").Value = 17630
$ws.Range("J89").Value = 17630
$ws.Range("L89").Value = 88150
$ws.Range("N89").Value = -99382
$ws.Range("H132").Value = 1537.2424
$ws.Range("I132").Value = 1211.05
$ws.Range("J132").Value = 2039.0769
$ws.Range("K132").Value = 3633.15
$ws.Range("L132").Value = 6117.2307
$ws.Range("M132").Value = -1103.15
$ws.Range("N132").Value = -11177.2307
